$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '22.001.90'
$ws.Cells.Item(2, 5).Value = '  -2.11%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.556.85'
$ws.Cells.Item(3, 5).Value = '  -1.17%  '
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  +0.34%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '1.003'
$ws.Cells.Item(5, 5).Value = '  +0.26%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '287.43'
$ws.Cells.Item(6, 5).Value = '  -0.25%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.3764'
$ws.Cells.Item(7, 5).Value = '  +1.70%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3252'
$ws.Cells.Item(8, 5).Value = '  -2.35%  '
$ws.Cells.Item(9, 2).Value = 'Polygon'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '1.127'
$ws.Cells.Item(9, 5).Value = '  -2.10%  '
$ws.Cells.Item(10, 2).Value = 'OKB'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '40.97'
$ws.Cells.Item(10, 5).Value = '  -14.30%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07303'
$ws.Cells.Item(11, 5).Value = '  -3.55%  '
$ws.Cells.Item(12, 5).Value = '  +0.36%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '19.72'
$ws.Cells.Item(13, 5).Value = '  -5.45%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '5.746'
$ws.Cells.Item(14, 5).Value = '  -3.60%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '6.843'
$ws.Cells.Item(15, 5).Value = '  -1.56%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '1.562.29'
$ws.Cells.Item(16, 5).Value = '  -0.33%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.00001083'
$ws.Cells.Item(17, 5).Value = '  -3.58%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.06632'
$ws.Cells.Item(18, 5).Value = '  -1.48%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '85.19'
$ws.Cells.Item(19, 5).Value = '  -3.54%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.423'
$ws.Cells.Item(20, 5).Value = '  +0.38%  '
$ws.Cells.Item(21, 5).Value = '  +0.17%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '15.97'
$ws.Cells.Item(22, 5).Value = '  -3.55%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '11.59'
$ws.Cells.Item(23, 5).Value = '  -3.83%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '22.013.78'
$ws.Cells.Item(24, 5).Value = '  -2.04%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.256'
$ws.Cells.Item(25, 5).Value = '  -5.48%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.531'
$ws.Cells.Item(26, 5).Value = '  -4.07%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '149.51'
$ws.Cells.Item(27, 5).Value = '  -1.17%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '18.92'
$ws.Cells.Item(28, 5).Value = '  -3.90%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '4.854'
$ws.Cells.Item(29, 5).Value = '  -2.81%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.739.68'
$ws.Cells.Item(30, 5).Value = '  -0.37%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '120.44'
$ws.Cells.Item(31, 5).Value = '  -4.06%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '1.116'
$ws.Cells.Item(32, 5).Value = '  +1.99%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '5.965'
$ws.Cells.Item(33, 5).Value = '  -2.53%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.763'
$ws.Cells.Item(34, 5).Value = '  -11.27%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '9.275'
$ws.Cells.Item(35, 5).Value = '  -6.21%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.08100'
$ws.Cells.Item(36, 5).Value = '  -3.11%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.219'
$ws.Cells.Item(37, 5).Value = '  -2.82%  '
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.02275'
$ws.Cells.Item(38, 5).Value = '  -7.71%  '
$ws.Cells.Item(39, 2).Value = 'Hedera'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.06129'
$ws.Cells.Item(39, 5).Value = '  -4.14%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.2124'
$ws.Cells.Item(40, 5).Value = '  -5.27%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.212'
$ws.Cells.Item(41, 5).Value = '  -6.70%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '10.88'
$ws.Cells.Item(42, 5).Value = '  -5.43%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.003'
$ws.Cells.Item(43, 5).Value = '  +0.26%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.5931'
$ws.Cells.Item(44, 5).Value = '  -5.75%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '13.47'
$ws.Cells.Item(45, 5).Value = '  -4.11%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '3.726'
$ws.Cells.Item(46, 5).Value = '  -1.46%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.5728'
$ws.Cells.Item(47, 5).Value = '  -6.47%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.951'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '119.89'
$ws.Cells.Item(49, 5).Value = '  -4.41%  '
$ws.Cells.Item(50, 5).Value = '  -4.65%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.06944'
$ws.Cells.Item(51, 5).Value = '  -3.78%  '
